$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "19.12.2023"
$ws.Range("D10").Value = "13:00-15:00; 16:00-20:00"
$ws.Range("C10").Value = "DB + request handling start"
$ws.Range("B10").Value = 360
$ws.Range("D10").NumberFormat = "h:mm"

$ws.Range("B11").Select()
